# Evan Schober Server Log.xlsx - Lab 8 update
# Adds 11 new log entries (rows 80-90) documenting Lab 8 AD delegation /
# group-management work performed on 3/2/2017 (serial date 42796), which
# previously occupied blank placeholder rows at the bottom of the log.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Date used for every new entry: 3/2/2017
$logDate = 42796

# row -> (Action, WhoAffected, RowHeight-or-$null)
$rows = @(
    @{ R = 80; Action = "Create Managers distribution group";                                         Who = "None";                                                       H = $null },
    @{ R = 81; Action = "Add Accounting managers to Managers group";                                   Who = "mbarnes";                                                    H = $null },
    @{ R = 82; Action = "Add Marketing managers to Managers group";                                    Who = "jsuarez";                                                    H = $null },
    @{ R = 83; Action = "Add Research-Dev managers to Managers group";                                 Who = "akimbly";                                                    H = 30    },
    @{ R = 84; Action = "Add Sales managers to Managers group";                                        Who = "mburnes";                                                    H = $null },
    @{ R = 85; Action = "Add Support managers to Managers group";                                      Who = "semery";                                                     H = $null },
    @{ R = 86; Action = "Create Support Resources group";                                              Who = "None";                                                       H = $null },
    @{ R = 87; Action = "Add Support group to Support Resources group";                                Who = "jrons, tplask, semery";                                      H = $null },
    @{ R = 88; Action = "Delegate control of domain PW Resets to PasswordAdmins group";                Who = "All";                                                        H = 30    },
    @{ R = 89; Action = "Delegate control of adding computers to the Domain to ComputerAdmins group";  Who = "All";                                                        H = 30    },
    @{ R = 90; Action = "Delegate control of departmental OUs to GPOLinkAdmins";                       Who = "Accounting,Marketing, Research-Dev, Sales, & Support OUs "; H = 45    }
)

# Column A on the existing rows above already carries the date-formatted
# (border + wrap + mm-dd-yy) style we need. Copy that formatting down onto
# the new rows instead of re-deriving it with NumberFormat, so the new
# cells share the same style record rather than each minting a new one.
$ws.Range("A79").Copy()

foreach ($row in $rows) {
    $r = $row.R

    $ws.Range("A$r").PasteSpecial(-4122)
    $ws.Range("A$r").Value = $logDate

    $ws.Range("B$r").Value = $row.Action
    $ws.Range("C$r").Value = "No"
    $ws.Range("D$r").Value = "N/A"
    $ws.Range("E$r").Value = $row.Who
    $ws.Range("F$r").Value = "Evan"
    $ws.Range("G$r").Value = "ES"

    if ($row.H) {
        $ws.Rows.Item($r).RowHeight = $row.H
    }
}

$excel.CutCopyMode = $false

# Update the saved view state to match where the editor last left off.
$ws.Activate()
$ws.Range("A73").Select()
$excel.ActiveWindow.ScrollRow = 73
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A91").Select()
